# Adds a new "enrollment_status" column to the sample_students sheet
# (inserted right before the existing "school_year" column) and backfills
# new mother/guardian name+contact data plus refreshed father_name /
# father_contact values and enrollment status per student row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column at T ("school_year" and everything after it
#    shifts one column to the right: T->U, U->V, V->W).
$ws.Columns("T:T").Insert()

# 2. New column header
$ws.Range("T1").Value = "enrollment_status"

# 3. Row 2 (John Doe Smith) - mother/guardian info + enrollment status
$ws.Range("O2").Value = "Mother 1"
$ws.Range("P2").Value = "Guardian 1"
$ws.Range("R2").Value = 9123456710
$ws.Range("S2").Value = 9123456722
$ws.Range("T2").Value = "Enrolled"

# 4. Row 3 (Jane Anne Doe) - father info updated + guardian/contact + status
$ws.Range("N3").Value = "Father 2"
$ws.Range("P3").Value = "Guardian 2"
$ws.Range("Q3").Value = 9123456781
$ws.Range("R3").Value = 9123456792
$ws.Range("S3").Value = 9123456733
$ws.Range("T3").Value = "Not Enrolled"
$ws.Range("U3").Value = "2023-2024"

# 5. Row 4 (Mark David Johnson) - father/mother info + contacts + status
$ws.Range("N4").Value = "Father 3"
$ws.Range("O4").Value = "Mother 3"
$ws.Range("Q4").Value = 9123456790
$ws.Range("R4").Value = 9123456709
$ws.Range("T4").Value = "Enrolled"

# 6. Row 5 (Sarah Lee Adams) - father info updated + guardian/contacts + status
$ws.Range("N5").Value = "Father 4"
$ws.Range("P5").Value = "Guardian 3"
$ws.Range("Q5").Value = 9123456708
$ws.Range("R5").Value = 9123456715
$ws.Range("S5").Value = 9123456734
$ws.Range("T5").Value = "Not Enrolled"

# 7. Row 6 (Alex James White) - father info updated + guardian/contact + status
$ws.Range("N6").Value = "Father 5"
$ws.Range("P6").Value = "Guardian 4"
$ws.Range("R6").Value = 9123456999
$ws.Range("T6").Value = "Graduate"

# 8. Row 7 (Natalie Rose Brown) - father info updated + guardian/contact + status
$ws.Range("N7").Value = "Father 6"
$ws.Range("P7").Value = "Guardian 5"
$ws.Range("R7").Value = 9123456877
$ws.Range("T7").Value = "Graduate"

# 9. Reset view back to the top-left corner of the sheet.
$ws.Range("A1").Select()
